## Update the Year-End Reconciliation template's school-year references
## from the 2019-2020 school year to the 2021-2022 school year.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Prior Year Reconciliation")

# Report title
$ws.Range("E2").Value = "RECONCILIATION REPORT FOR THE 2021-2022 SCHOOL YEAR"

# Column header "Total Amount Due for ... School Year"
$ws.Range("H9").Value = "Total Amount Due for 2021-2022 School Year"

# Leading-space / quote-prefixed label (preserve the quote-prefix style by
# typing a leading apostrophe, just like a user would in Excel)
$ws.Range("G12").Value = "'          Total Amount Due for 2021-2022 School Year:"

# Month/year labels in the payment schedule
$ws.Range("C16").Value = "July, 2021"
$ws.Range("C22").Value = "January, 2022"

# Leading-space / quote-prefixed label (preserve the quote-prefix style by
# typing a leading apostrophe, just like a user would in Excel)
$ws.Range("G30").Value = "'            Total Paid to Date for 2021-2022 School Year:"

# Footnote about enrollment cut-off date
$ws.Range("A34").Value = "1.  Do not include students that enrolled after`n     March 13, 2022."
